$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Ranger URL text (shared string) used by C2:C4 and by the
# hyperlinks that decorate them: "http://dev.ranger9dot5.xpms.ai" ->
# "http://dev.ranger.xpms.ai"
$newUrl = "http://dev.ranger.xpms.ai"

$ws.Range("C2").Value = $newUrl
$ws.Range("C3").Value = $newUrl
$ws.Range("C4").Value = $newUrl

# Refresh the hyperlinks backing those cells so both the display text and
# the link target point at the new address.
$ws.Hyperlinks.Delete()

$h2 = $ws.Range("C2").Hyperlinks.Item(1)
$h2.Address = $newUrl
$h2.TextToDisplay = $newUrl

$h3 = $ws.Range("C3").Hyperlinks.Item(1)
$h3.Address = $newUrl
$h3.TextToDisplay = $newUrl

$h4 = $ws.Range("C4").Hyperlinks.Item(1)
$h4.Address = $newUrl
$h4.TextToDisplay = $newUrl

# Move the active selection from B2 to C4.
$ws.Range("C4").Select() | Out-Null
